# Insert a new "2022-Q1" sheet after "2021-Q4" (before "总计"),
# populate it with per-fund holding data, and prepend a matching summary
# row to the "总计" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Create the new sheet right after "2021-Q4" (i.e. right before "总计")
# ---------------------------------------------------------------------
$q4Sheet = $wb.Worksheets.Item("2021-Q4")
$newSheet = $wb.Worksheets.Add($null, $q4Sheet)
$newSheet.Name = "2022-Q1"

# ---------------------------------------------------------------------
# 2) Header row (bold, bordered, centered-top — matches the other
#    per-quarter sheets such as "2021-Q4")
# ---------------------------------------------------------------------
$hdr = $newSheet.Range("B1:H1")
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160
$hdr.Borders.LineStyle = 1

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# ---------------------------------------------------------------------
# 3) Per-fund holding rows
# ---------------------------------------------------------------------
$fundRows = @(
    @{B="004224"; C="南方军工改革灵活配置混合A"; D="63.35"; E="86.90"; F="5.57"; G="3.5286"; H=10},
    @{B="001513"; C="易方达信息产业混合"; D="32.50"; E="92.37"; F="4.18"; G="1.3585"; H=1},
    @{B="506002"; C="易方达科创板两年定期开放混合型证券投资基金"; D="29.21"; E="86.29"; F="3.22"; G="0.9406"; H=8},
    @{B="011148"; C="南方军工改革灵活配置混合C"; D="15.56"; E="86.90"; F="5.57"; G="0.8667"; H=10},
    @{B="010410"; C="长城品质成长混合A"; D="29.08"; E="70.65"; F="2.96"; G="0.8608"; H=3},
    @{B="506003"; C="富国科创板两年定期开放混合"; D="23.60"; E="98.48"; F="3.18"; G="0.7505"; H=3},
    @{B="506000"; C="南方科创板 3 年定期开放混合"; D="24.62"; E="96.87"; F="2.91"; G="0.7164"; H=9},
    @{B="160642"; C="鹏华增瑞灵活配置混合(LOF)"; D="6.76"; E="91.34"; F="10.01"; G="0.6767"; H=1},
    @{B="110002"; C="易方达策略成长混合"; D="12.15"; E="88.76"; F="4.93"; G="0.5990"; H=2},
    @{B="112002"; C="易方达策略成长二号混合"; D="10.72"; E="87.99"; F="5.22"; G="0.5596"; H=1},
    @{B="200012"; C="长城中小盘成长混合"; D="12.65"; E="84.26"; F="3.70"; G="0.4680"; H=2},
    @{B="200007"; C="长城安心回报混合"; D="11.42"; E="71.80"; F="2.97"; G="0.3392"; H=4},
    @{B="001076"; C="易方达改革红利混合"; D="8.95"; E="88.46"; F="3.49"; G="0.3124"; H=10},
    @{B="005738"; C="长城智能产业灵活配置混合"; D="8.04"; E="84.35"; F="3.85"; G="0.3095"; H=2},
    @{B="010284"; C="长城价值成长六个月持有期混合A"; D="7.14"; E="87.09"; F="3.13"; G="0.2235"; H=2},
    @{B="006769"; C="长城研究精选混合"; D="6.41"; E="82.49"; F="3.09"; G="0.1981"; H=2},
    @{B="010602"; C="长城均衡优选混合"; D="4.55"; E="84.61"; F="3.45"; G="0.1570"; H=2},
    @{B="010824"; C="天弘创新成长混合A"; D="3.54"; E="82.73"; F="3.90"; G="0.1381"; H=9},
    @{B="005310"; C="广发电子信息传媒产业精选股票A"; D="3.99"; E="90.16"; F="3.37"; G="0.1345"; H=9},
    @{B="506008"; C="长城科创两年定开混合A"; D="3.57"; E="62.09"; F="3.04"; G="0.1085"; H=3},
    @{B="005495"; C="创金合信科技成长主题股票A"; D="2.62"; E="84.91"; F="3.83"; G="0.1003"; H=2},
    @{B="010411"; C="长城品质成长混合C"; D="2.34"; E="70.65"; F="2.96"; G="0.0693"; H=3},
    @{B="010495"; C="创金合信创新驱动股票A"; D="1.47"; E="82.29"; F="3.71"; G="0.0545"; H=6},
    @{B="010825"; C="天弘创新成长混合C"; D="1.08"; E="82.73"; F="3.90"; G="0.0421"; H=9},
    @{B="005496"; C="创金合信科技成长主题股票C"; D="0.72"; E="84.91"; F="3.83"; G="0.0276"; H=2},
    @{B="010236"; C="广发电子信息传媒产业精选股票C"; D="0.81"; E="90.16"; F="3.37"; G="0.0273"; H=9},
    @{B="010285"; C="长城价值成长六个月持有期混合C"; D="0.52"; E="87.09"; F="3.13"; G="0.0163"; H=2},
    @{B="010496"; C="创金合信创新驱动股票C"; D="0.37"; E="82.29"; F="3.71"; G="0.0137"; H=6},
    @{B="004521"; C="安信工业4.0主题沪港深精选灵活配置混合A"; D="0.09"; E="85.98"; F="4.76"; G="0.0043"; H=9},
    @{B="012793"; C="长城科创两年定开混合C"; D="0.12"; E="62.09"; F="3.04"; G="0.0036"; H=3},
    @{B="004522"; C="安信工业4.0主题沪港深精选灵活配置混合C"; D="0.05"; E="85.98"; F="4.76"; G="0.0024"; H=9}
)

$r = 2
foreach ($fr in $fundRows) {
    $idxCell = $newSheet.Range("A$r")
    $idxCell.Value = $r - 2
    $idxCell.Font.Bold = $true
    $idxCell.HorizontalAlignment = -4108
    $idxCell.VerticalAlignment = -4160
    $idxCell.Borders.LineStyle = 1

    # Text columns (B..G) must stay text, even though several look numeric
    # (leading zeros in fund codes, fixed-decimal percentages, etc.), so
    # force a text number format before assigning, then drop back to the
    # Normal cell style to avoid leaving a stray quote-prefix/text format.
    $txt = $newSheet.Range("B$r`:G$r")
    $txt.NumberFormat = "@"
    $newSheet.Range("B$r").Value = $fr.B
    $newSheet.Range("C$r").Value = $fr.C
    $newSheet.Range("D$r").Value = $fr.D
    $newSheet.Range("E$r").Value = $fr.E
    $newSheet.Range("F$r").Value = $fr.F
    $newSheet.Range("G$r").Value = $fr.G
    $txt.Style = "Normal"

    $newSheet.Range("H$r").Value = $fr.H

    $r++
}

$null = $newSheet.Range("A1").Select()

# ---------------------------------------------------------------------
# 4) Prepend the 2022-Q1 roll-up row to "总计", shifting the existing
#    rows down by exactly one and renumbering the index column (A) so
#    it stays a contiguous 0-based sequence.
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$null = $totalSheet.Range("A2:D2").Insert()

$idxCell = $totalSheet.Range("A2")
$idxCell.Value = 0
$idxCell.Font.Bold = $true
$idxCell.HorizontalAlignment = -4108
$idxCell.VerticalAlignment = -4160
$idxCell.Borders.LineStyle = 1

$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 31
$totalSheet.Range("D2").Value = 13.61

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
$totalSheet.Range("A6").Value = 4
$totalSheet.Range("A7").Value = 5

$null = $totalSheet.Range("A1").Select()
